$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.591135333333333
$ws.Range("H2").Value = 4.773406
$ws.Range("I2").Value = 0.4707829234247397
$ws.Range("J2").Value = 0.4707829234247397
$ws.Range("M2").Value = 14.25737566666667
$ws.Range("N2").Value = 42.772127
$ws.Range("O2").Value = 0.2087950866344732
$ws.Range("P2").Value = 0.2087950866344732
$ws.Range("Q2").Value = 22.68541418384022
$ws.Range("R2").Value = 204.168727654562
$ws.Range("S2").Value = 0.09829716128249907
$ws.Range("T2").Value = 0.0982971612824991

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.591135333333333
$ws.Range("H3").Value = 4.773406
$ws.Range("I3").Value = 0.4707829234247397
$ws.Range("J3").Value = 0.4707829234247397
$ws.Range("N3").Value = 87.128332
$ws.Range("O3").Value = 0.4253229592313036
$ws.Range("P3").Value = 0.4253229592313036
$ws.Range("Q3").Value = 46.2109891931991
$ws.Range("R3").Value = 415.8989027387919
$ws.Range("S3").Value = 0.2002347861465745
$ws.Range("T3").Value = 0.2002347861465745

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.591135333333333
$ws.Range("H4").Value = 4.773406
$ws.Range("I4").Value = 0.4707829234247397
$ws.Range("J4").Value = 0.4707829234247397
$ws.Range("M4").Value = 20.11084633333333
$ws.Range("N4").Value = 60.332539
$ws.Range("O4").Value = 0.2945174484164121
$ws.Range("P4").Value = 0.2945174484164122
$ws.Range("Q4").Value = 31.99907818420377
$ws.Range("R4").Value = 287.991703657834
$ws.Range("S4").Value = 0.1386537853650735
$ws.Range("T4").Value = 0.1386537853650735

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.591135333333333
$ws.Range("H5").Value = 4.773406
$ws.Range("I5").Value = 0.4707829234247397
$ws.Range("J5").Value = 0.4707829234247397
$ws.Range("M5").Value = 4.873057999999999
$ws.Range("N5").Value = 14.619174
$ws.Range("O5").Value = 0.07136450571781097
$ws.Range("P5").Value = 0.07136450571781099
$ws.Range("Q5").Value = 7.753694765182665
$ws.Range("R5").Value = 69.78325288664399
$ws.Range("S5").Value = 0.0335971906305926
$ws.Range("T5").Value = 0.03359719063059261

# Row 6
$ws.Range("I6").Value = 0.3035973020998604
$ws.Range("J6").Value = 0.3035973020998604
$ws.Range("M6").Value = 14.25737566666667
$ws.Range("N6").Value = 42.772127
$ws.Range("O6").Value = 0.2087950866344732
$ws.Range("P6").Value = 0.2087950866344732
$ws.Range("Q6").Value = 14.62931257814156
$ws.Range("R6").Value = 131.663813203274
$ws.Range("S6").Value = 0.06338962499393268
$ws.Range("T6").Value = 0.06338962499393269

# Row 7
$ws.Range("I7").Value = 0.3035973020998604
$ws.Range("J7").Value = 0.3035973020998604
$ws.Range("N7").Value = 87.128332
$ws.Range("O7").Value = 0.4253229592313036
$ws.Range("P7").Value = 0.4253229592313036
$ws.Range("S7").Value = 0.1291269029437527
$ws.Range("T7").Value = 0.1291269029437527

# Row 8
$ws.Range("I8").Value = 0.3035973020998604
$ws.Range("J8").Value = 0.3035973020998604
$ws.Range("M8").Value = 20.11084633333333
$ws.Range("N8").Value = 60.332539
$ws.Range("O8").Value = 0.2945174484164121
$ws.Range("P8").Value = 0.2945174484164122
$ws.Range("Q8").Value = 20.63548468524644
$ws.Range("R8").Value = 185.719362167218
$ws.Range("S8").Value = 0.08941470276055753
$ws.Range("T8").Value = 0.08941470276055755

# Row 9
$ws.Range("I9").Value = 0.3035973020998604
$ws.Range("J9").Value = 0.3035973020998604
$ws.Range("M9").Value = 4.873057999999999
$ws.Range("N9").Value = 14.619174
$ws.Range("O9").Value = 0.07136450571781097
$ws.Range("P9").Value = 0.07136450571781099
$ws.Range("Q9").Value = 5.000183088398666
$ws.Range("R9").Value = 45.001647795588
$ws.Range("S9").Value = 0.02166607140161747
$ws.Range("T9").Value = 0.02166607140161748

# Row 10
$ws.Range("G10").Value = 0.730693
$ws.Range("H10").Value = 2.192079
$ws.Range("I10").Value = 0.2161964349979826
$ws.Range("J10").Value = 0.2161964349979826
$ws.Range("M10").Value = 14.25737566666667
$ws.Range("N10").Value = 42.772127
$ws.Range("O10").Value = 0.2087950866344732
$ws.Range("P10").Value = 0.2087950866344732
$ws.Range("Q10").Value = 10.41776459800367
$ws.Range("R10").Value = 93.75988138203302
$ws.Range("S10").Value = 0.04514075337546803
$ws.Range("T10").Value = 0.04514075337546804

# Row 11
$ws.Range("G11").Value = 0.730693
$ws.Range("H11").Value = 2.192079
$ws.Range("I11").Value = 0.2161964349979826
$ws.Range("J11").Value = 0.2161964349979826
$ws.Range("N11").Value = 87.128332
$ws.Range("O11").Value = 0.4253229592313036
$ws.Range("P11").Value = 0.4253229592313036
$ws.Range("Q11").Value = 21.22135409802533
$ws.Range("R11").Value = 190.992186882228
$ws.Range("S11").Value = 0.09195330750860013
$ws.Range("T11").Value = 0.09195330750860015

# Row 12
$ws.Range("G12").Value = 0.730693
$ws.Range("H12").Value = 2.192079
$ws.Range("I12").Value = 0.2161964349979826
$ws.Range("J12").Value = 0.2161964349979826
$ws.Range("M12").Value = 20.11084633333333
$ws.Range("N12").Value = 60.332539
$ws.Range("O12").Value = 0.2945174484164121
$ws.Range("P12").Value = 0.2945174484164122
$ws.Range("Q12").Value = 14.69485463984233
$ws.Range("R12").Value = 132.253691758581
$ws.Range("S12").Value = 0.06367362239233054
$ws.Range("T12").Value = 0.06367362239233056

# Row 13
$ws.Range("G13").Value = 0.730693
$ws.Range("H13").Value = 2.192079
$ws.Range("I13").Value = 0.2161964349979826
$ws.Range("J13").Value = 0.2161964349979826
$ws.Range("M13").Value = 4.873057999999999
$ws.Range("N13").Value = 14.619174
$ws.Range("O13").Value = 0.07136450571781097
$ws.Range("P13").Value = 0.07136450571781099
$ws.Range("Q13").Value = 3.560709369194
$ws.Range("R13").Value = 32.046384322746
$ws.Range("S13").Value = 0.01542875172158388
$ws.Range("T13").Value = 0.01542875172158388

# Row 14
$ws.Range("G14").Value = 0.03184866666666667
$ws.Range("H14").Value = 0.09554600000000001
$ws.Range("I14").Value = 0.009423339477417213
$ws.Range("J14").Value = 0.009423339477417213
$ws.Range("M14").Value = 14.25737566666667
$ws.Range("N14").Value = 42.772127
$ws.Range("O14").Value = 0.2087950866344732
$ws.Range("P14").Value = 0.2087950866344732
$ws.Range("Q14").Value = 0.4540784051491112
$ws.Range("R14").Value = 4.086705646342001
$ws.Range("S14").Value = 0.001967546982573378
$ws.Range("T14").Value = 0.001967546982573378

# Row 15
$ws.Range("G15").Value = 0.03184866666666667
$ws.Range("H15").Value = 0.09554600000000001
$ws.Range("I15").Value = 0.009423339477417213
$ws.Range("J15").Value = 0.009423339477417213
$ws.Range("N15").Value = 87.128332
$ws.Range("O15").Value = 0.4253229592313036
$ws.Range("P15").Value = 0.4253229592313036
$ws.Range("Q15").Value = 0.9249737343635557
$ws.Range("R15").Value = 8.324763609272001
$ws.Range("S15").Value = 0.004007962632376255
$ws.Range("T15").Value = 0.004007962632376255

# Row 16
$ws.Range("G16").Value = 0.03184866666666667
$ws.Range("H16").Value = 0.09554600000000001
$ws.Range("I16").Value = 0.009423339477417213
$ws.Range("J16").Value = 0.009423339477417213
$ws.Range("M16").Value = 20.11084633333333
$ws.Range("N16").Value = 60.332539
$ws.Range("O16").Value = 0.2945174484164121
$ws.Range("P16").Value = 0.2945174484164122
$ws.Range("Q16").Value = 0.6405036412548889
$ws.Range("R16").Value = 5.764532771294
$ws.Range("S16").Value = 0.002775337898450564
$ws.Range("T16").Value = 0.002775337898450564

# Row 17
$ws.Range("G17").Value = 0.03184866666666667
$ws.Range("H17").Value = 0.09554600000000001
$ws.Range("I17").Value = 0.009423339477417213
$ws.Range("J17").Value = 0.009423339477417213
$ws.Range("M17").Value = 4.873057999999999
$ws.Range("N17").Value = 14.619174
$ws.Range("O17").Value = 0.07136450571781097
$ws.Range("P17").Value = 0.07136450571781099
$ws.Range("Q17").Value = 0.1552003998893333
$ws.Range("R17").Value = 1.396803599004
$ws.Range("S17").Value = 0.0006724919640170145
$ws.Range("T17").Value = 0.0006724919640170146
